# Update the cryptocurrency Price (col D) and Volume(1h) (col E) columns
# on the active sheet to the latest scraped values (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Cells.Item(2, 4).Value = '29.453.55'
$ws.Cells.Item(2, 5).Value = '  +0.39%  '

# Row 3: Ethereum
$ws.Cells.Item(3, 4).Value = '1.850.87'
$ws.Cells.Item(3, 5).Value = '  +0.45%  '

# Row 4: TetherUSD
$ws.Cells.Item(4, 4).Value = '''1.0000'
$ws.Cells.Item(4, 5).Value = '  +0.13%  '

# Row 5: BNB
$ws.Cells.Item(5, 4).Value = '''240.98'
$ws.Cells.Item(5, 5).Value = '  +0.78%  '

# Row 6: XRP
$ws.Cells.Item(6, 4).Value = '''0.6296'
$ws.Cells.Item(6, 5).Value = '  +0.01%  '

# Row 7: USDC
$ws.Cells.Item(7, 5).Value = '  +0.09%  '

# Row 8: Dogecoin
$ws.Cells.Item(8, 4).Value = '''0.07701'
$ws.Cells.Item(8, 5).Value = '  +2.27%  '

# Row 9: Cardano
$ws.Cells.Item(9, 4).Value = '''0.2929'
$ws.Cells.Item(9, 5).Value = '  -0.54%  '

# Row 10: Solana
$ws.Cells.Item(10, 4).Value = '''24.72'
$ws.Cells.Item(10, 5).Value = '  +0.91%  '

# Row 11: TRON
$ws.Cells.Item(11, 4).Value = '''0.07744'
$ws.Cells.Item(11, 5).Value = '  +0.73%  '

# Row 12: WrappedEther
$ws.Cells.Item(12, 4).Value = '1.885.96'
$ws.Cells.Item(12, 5).Value = '  +1.16%  '

# Row 13: Polkadot
$ws.Cells.Item(13, 5).Value = '  +1.21%  '

# Row 14: ShibaInu
$ws.Cells.Item(14, 4).Value = '''0.00001078'
$ws.Cells.Item(14, 5).Value = '  +4.40%  '

# Row 15: Polygon
$ws.Cells.Item(15, 4).Value = '''0.6792'
$ws.Cells.Item(15, 5).Value = '  +0.17%  '

# Row 16: Litecoin
$ws.Cells.Item(16, 4).Value = '''83.70'
$ws.Cells.Item(16, 5).Value = '  +0.76%  '

# Row 17: WrappedliquidstakedEther2.0
$ws.Cells.Item(17, 4).Value = '2.156.83'
$ws.Cells.Item(17, 5).Value = '  +2.11%  '

# Row 19: WrappedBTC
$ws.Cells.Item(19, 4).Value = '29.496.75'
$ws.Cells.Item(19, 5).Value = '  +0.43%  '

# Row 20: BitcoinCash
$ws.Cells.Item(20, 4).Value = '''228.54'
$ws.Cells.Item(20, 5).Value = '  -0.11%  '

# Row 21: Avalanche
$ws.Cells.Item(21, 4).Value = '''12.45'
$ws.Cells.Item(21, 5).Value = '  +0.32%  '

# Row 22: Dai
$ws.Cells.Item(22, 5).Value = '  +0.02%  '

# Row 23: Chainlink
$ws.Cells.Item(23, 4).Value = '''7.448'
$ws.Cells.Item(23, 5).Value = '  -0.08%  '

# Row 24: BinanceUSD
$ws.Cells.Item(24, 4).Value = '''1.001'
$ws.Cells.Item(24, 5).Value = '  +0.10%  '

# Row 25: Monero
$ws.Cells.Item(25, 4).Value = '''157.57'
$ws.Cells.Item(25, 5).Value = '  +0.82%  '

# Row 26: Stellar
$ws.Cells.Item(26, 4).Value = '''0.1378'
$ws.Cells.Item(26, 5).Value = '  -1.09%  '

# Row 27: Cosmos
$ws.Cells.Item(27, 4).Value = '''8.410'
$ws.Cells.Item(27, 5).Value = '  +0.68%  '

# Row 28: EthereumClassic
$ws.Cells.Item(28, 5).Value = '  +0.55%  '

# Row 29: Toncoin
$ws.Cells.Item(29, 4).Value = '''1.345'
$ws.Cells.Item(29, 5).Value = '  +6.05%  '

# Row 30: PancakeSwap
$ws.Cells.Item(30, 4).Value = '''1.469'
$ws.Cells.Item(30, 5).Value = '  +0.59%  '

# Row 31: Hedera
$ws.Cells.Item(31, 5).Value = '  +0.63%  '

# Row 32: Filecoin
$ws.Cells.Item(32, 5).Value = '  +0.46%  '

# Row 33: InternetComputer(DFINITY)
$ws.Cells.Item(33, 4).Value = '''4.036'
$ws.Cells.Item(33, 5).Value = '  +0.42%  '

# Row 34: LidoDAOToken
$ws.Cells.Item(34, 4).Value = '''1.847'
$ws.Cells.Item(34, 5).Value = '  +1.04%  '

# Row 35: ARBITRUM
$ws.Cells.Item(35, 5).Value = '  +0.85%  '

# Row 36: ImmutableX
$ws.Cells.Item(36, 4).Value = '''0.7033'
$ws.Cells.Item(36, 5).Value = '  -0.90%  '

# Row 37: HuobiToken
$ws.Cells.Item(37, 4).Value = '''2.585'
$ws.Cells.Item(37, 5).Value = '  -0.14%  '

# Row 38: MXToken
$ws.Cells.Item(38, 4).Value = '''2.782'
$ws.Cells.Item(38, 5).Value = '  +0.44%  '

# Row 39: VeChain
$ws.Cells.Item(39, 4).Value = '''0.01793'
$ws.Cells.Item(39, 5).Value = '  -0.78%  '

# Row 40: Maker
$ws.Cells.Item(40, 4).Value = '1.220.23'
$ws.Cells.Item(40, 5).Value = '  -1.68%  '

# Row 41: FraxShare
$ws.Cells.Item(41, 4).Value = '''6.551'
$ws.Cells.Item(41, 5).Value = '  +5.06%  '

# Row 42: TrustWalletToken
$ws.Cells.Item(42, 4).Value = '''0.9065'
$ws.Cells.Item(42, 5).Value = '  +0.61%  '

# Row 43: PaxDollar
$ws.Cells.Item(43, 5).Value = '  +0.16%  '

# Row 44: Quant
$ws.Cells.Item(44, 5).Value = '  +0.15%  '

# Row 45: Aave
$ws.Cells.Item(45, 4).Value = '''66.32'
$ws.Cells.Item(45, 5).Value = '  +1.25%  '

# Row 46: BabyDogeCoin
$ws.Cells.Item(46, 4).Value = '''0.00000000119'
$ws.Cells.Item(46, 5).Value = '  +1.52%  '

# Row 47: Aptos
$ws.Cells.Item(47, 4).Value = '''7.145'
$ws.Cells.Item(47, 5).Value = '  +0.65%  '

# Row 48: TheSandbox
$ws.Cells.Item(48, 4).Value = '''0.4020'
$ws.Cells.Item(48, 5).Value = '  +0.68%  '

# Row 49: EnergySwap
$ws.Cells.Item(49, 4).Value = '''9.033'
$ws.Cells.Item(49, 5).Value = '  +1.08%  '

# Row 50: RenderToken
$ws.Cells.Item(50, 4).Value = '''1.682'
$ws.Cells.Item(50, 5).Value = '  +0.45%  '

# Row 51: Algorand
$ws.Cells.Item(51, 5).Value = '  +2.40%  '
